# Add the new "EmployeeListPage" worksheet after "AddEmployeePage" and
# populate it with the Employee List test data, then update the active
# sheet / selection state to match.

$wb = $excel.ActiveWorkbook
$addEmployeeSheet = $wb.Worksheets.Item("AddEmployeePage")

# --- Add the new worksheet at the end, named EmployeeListPage ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "EmployeeListPage"

# --- Reuse the existing header / body cell formatting from AddEmployeePage
#     (style "left/top aligned bold" for the header row, "left/top aligned"
#     for the body rows) instead of inventing new styles. ---
$addEmployeeSheet.Range("A1").Copy() | Out-Null
$newSheet.Range("A1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$addEmployeeSheet.Range("A2").Copy() | Out-Null
$newSheet.Range("A2:A8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Header row ---
$newSheet.Range("A1").Value = "TestCaseName"
$newSheet.Range("B1").Value = "EmployeeID"
$newSheet.Range("C1").Value = "EmployeeName"
$newSheet.Range("D1").Value = "SupervisorName"
$newSheet.Range("E1").Value = "Employment Status"
$newSheet.Range("F1").Value = "Include"
$newSheet.Range("G1").Value = "JobTitle"
$newSheet.Range("H1").Value = "SubUnit"
$newSheet.Range("I1").Value = "Expected"
$newSheet.Range("J1").Value = "Actual"

# --- TestCaseName column ---
$newSheet.Range("A2").Value = "testCase01"
$newSheet.Range("A3").Value = "testCase02"
$newSheet.Range("A4").Value = "testCase03"
$newSheet.Range("A5").Value = "testCase04"
$newSheet.Range("A6").Value = "testCase05"
$newSheet.Range("A7").Value = "testCase06"
$newSheet.Range("A8").Value = "testCase07"

# --- Previously-active sheet (AddEmployeePage) selection changes ---
$addEmployeeSheet.Range("G1:H1").Select() | Out-Null

# --- Selection / active cell + activation on new sheet (now the active tab) ---
$newSheet.Activate()
$newSheet.Range("I5").Select() | Out-Null
